# Scheduled-runner update: refresh Universalis market-price snapshots and the
# derived Leve profit columns (H:N) for the affected leves on each crafting-class
# sheet. A few rows also pick up/lose a trailing HQ-profit cell as a side effect
# of the refreshed source data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1334
$ws.Range("J6").Value = 636.4
$ws.Range("L6").Value = 1909.2
$ws.Range("N6").Value = -2133.2
$ws.Range("H8").Value = 233.66667
$ws.Range("I8").Value = 100.5
$ws.Range("K8").Value = 301.5
$ws.Range("M8").Value = -162.5
$ws.Range("H40").Value = 3626.2666
$ws.Range("I40").Value = 2815
$ws.Range("J40").Value = 3829.0833
$ws.Range("K40").Value = 2815
$ws.Range("L40").Value = 3829.0833
$ws.Range("M40").Value = -2640
$ws.Range("N40").Value = -4179.0833
$ws.Range("H116").Value = 8037.8823
$ws.Range("I116").Value = 7973.9165
$ws.Range("J116").Value = 8191.4
$ws.Range("K116").Value = 7973.9165
$ws.Range("L116").Value = 8191.4
$ws.Range("M116").Value = -4531.9165
$ws.Range("N116").Value = -15075.4
$ws.Range("H129").Value = 1503.0714
$ws.Range("I129").Value = 840.5
$ws.Range("K129").Value = 2521.5
$ws.Range("M129").Value = 2478.5
$ws.Range("H137").Value = 2416698.5
$ws.Range("I137").Value = 908.25
$ws.Range("J137").Value = 3705120
$ws.Range("K137").Value = 2724.75
$ws.Range("L137").Value = 11115360
$ws.Range("M137").Value = -174.75
$ws.Range("N137").Value = -11120460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1925.963
$ws.Range("I2").Value = 1585.5883
$ws.Range("J2").Value = 2504.6
$ws.Range("K2").Value = 1585.5883
$ws.Range("L2").Value = 2504.6
$ws.Range("M2").Value = -1472.5883
$ws.Range("N2").Value = -2730.6
$ws.Range("H32").Value = 28760820
$ws.Range("I32").Value = 35963836
$ws.Range("J32").Value = 6496952.5
$ws.Range("K32").Value = 35963836
$ws.Range("L32").Value = 6496952.5
$ws.Range("M32").Value = -35963549
$ws.Range("N32").Value = -6497526.5
$ws.Range("H97").Value = 3333.8462
$ws.Range("I97").Value = 1838.3334
$ws.Range("K97").Value = 1838.3334
$ws.Range("M97").Value = -1342.3334
$ws.Range("H116").Value = 1925.963
$ws.Range("I116").Value = 1585.5883
$ws.Range("J116").Value = 2504.6
$ws.Range("K116").Value = 1585.5883
$ws.Range("L116").Value = 2504.6
$ws.Range("M116").Value = 708.4117000000001
$ws.Range("N116").Value = -7092.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1925.963
$ws.Range("I3").Value = 1585.5883
$ws.Range("J3").Value = 2504.6
$ws.Range("K3").Value = 1585.5883
$ws.Range("L3").Value = 2504.6
$ws.Range("M3").Value = -1471.5883
$ws.Range("N3").Value = -2732.6
$ws.Range("H134").Value = 5105869.5
$ws.Range("I134").Value = 7145117.5
$ws.Range("J134").Value = 7749.75
$ws.Range("K134").Value = 21435352.5
$ws.Range("L134").Value = 23249.25
$ws.Range("M134").Value = -21432817.5
$ws.Range("N134").Value = -28319.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1962.125
$ws.Range("I16").Value = 2099.5715
$ws.Range("K16").Value = 2099.5715
$ws.Range("M16").Value = -1812.5715
$ws.Range("H31").Value = 3599.5435
$ws.Range("I31").Value = 1171.5834
$ws.Range("J31").Value = 4456.4707
$ws.Range("K31").Value = 1171.5834
$ws.Range("L31").Value = 4456.4707
$ws.Range("M31").Value = -876.5834
$ws.Range("N31").Value = -5046.4707
$ws.Range("H34").Value = 3599.5435
$ws.Range("I34").Value = 1171.5834
$ws.Range("J34").Value = 4456.4707
$ws.Range("K34").Value = 1171.5834
$ws.Range("L34").Value = 4456.4707
$ws.Range("M34").Value = -969.5834
$ws.Range("N34").Value = -4860.4707
$ws.Range("H53").Value = 107999
$ws.Range("J53").Value = 107999
$ws.Range("L53").Value = 107999
$ws.Range("N53").Value = -109213
$ws.Range("H58").Value = 2833.5112
$ws.Range("I58").Value = 2537.9487
$ws.Range("K58").Value = 2537.9487
$ws.Range("M58").Value = -2334.9487
$ws.Range("H99").Value = 2256
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 1614.5769
$ws.Range("I105").Value = 856.6667
$ws.Range("J105").Value = 2264.2144
$ws.Range("K105").Value = 856.6667
$ws.Range("L105").Value = 2264.2144
$ws.Range("M105").Value = 890.3333
$ws.Range("N105").Value = -5758.2144
$ws.Range("H113").Value = 1962.125
$ws.Range("I113").Value = 2099.5715
$ws.Range("K113").Value = 2099.5715
$ws.Range("M113").Value = 70.42849999999999
$ws.Range("H126").Value = 2256
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 2833.5112
$ws.Range("I136").Value = 2537.9487
$ws.Range("K136").Value = 7613.8461
$ws.Range("M136").Value = -5063.8461

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1946.2307
$ws.Range("I5").Value = 1759.6
$ws.Range("J5").Value = 2062.875
$ws.Range("K5").Value = 5278.799999999999
$ws.Range("L5").Value = 6188.625
$ws.Range("M5").Value = -5166.799999999999
$ws.Range("N5").Value = -6412.625
$ws.Range("H26").Value = 233.33333
$ws.Range("I26").Value = 200
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 600
$ws.Range("L26").Value = 900
$ws.Range("M26").Value = -312
$ws.Range("N26").Value = -1476
$ws.Range("H38").Value = 46.25
$ws.Range("I38").Value = 46.909092
$ws.Range("J38").Value = 45.444443
$ws.Range("K38").Value = 140.727276
$ws.Range("L38").Value = 136.333329
$ws.Range("M38").Value = 206.272724
$ws.Range("N38").Value = -830.333329
$ws.Range("H129").Value = 1895.5
$ws.Range("J129").Value = 2992.5
$ws.Range("L129").Value = 8977.5
$ws.Range("N129").Value = -18977.5
$ws.Range("H131").Value = 1717.3948
$ws.Range("J131").Value = 1783.0312
$ws.Range("L131").Value = 5349.0936
$ws.Range("N131").Value = -15429.0936
$ws.Range("H132").Value = 528866
$ws.Range("J132").Value = 1114465.4
$ws.Range("L132").Value = 10030188.6
$ws.Range("N132").Value = -10035248.6
$ws.Range("H135").Value = 1946.2307
$ws.Range("I135").Value = 1759.6
$ws.Range("J135").Value = 2062.875
$ws.Range("K135").Value = 15836.4
$ws.Range("L135").Value = 18565.875
$ws.Range("M135").Value = -13301.4
$ws.Range("N135").Value = -23635.875
$ws.Range("H137").Value = 7971
$ws.Range("I137").Value = 1750.5834
$ws.Range("J137").Value = 22900
$ws.Range("K137").Value = 5251.7502
$ws.Range("L137").Value = 68700
$ws.Range("M137").Value = -151.7502000000004
$ws.Range("N137").Value = -78900
$ws.Range("H138").Value = 16386221
$ws.Range("I138").Value = 1586.8
$ws.Range("J138").Value = 71001660
$ws.Range("K138").Value = 4760.4
$ws.Range("L138").Value = 213004980
$ws.Range("M138").Value = 379.6000000000004
$ws.Range("N138").Value = -213015260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3644
$ws.Range("I80").Value = 3608.3
$ws.Range("K80").Value = 3608.3
$ws.Range("M80").Value = -2610.3
$ws.Range("H83").Value = 3644
$ws.Range("I83").Value = 3608.3
$ws.Range("K83").Value = 18041.5
$ws.Range("M83").Value = -13049.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1080001.8
$ws.Range("I2").Value = 7500000
$ws.Range("K2").Value = 7500000
$ws.Range("M2").Value = -7499888
$ws.Range("H7").Value = 6518.5
$ws.Range("I7").Value = 6577.5454
$ws.Range("J7").Value = 6425.7144
$ws.Range("K7").Value = 6577.5454
$ws.Range("L7").Value = 6425.7144
$ws.Range("M7").Value = -6465.5454
$ws.Range("N7").Value = -6649.7144
$ws.Range("H22").Value = 4417.3335
$ws.Range("I22").Value = 3001
$ws.Range("J22").Value = 7250
$ws.Range("K22").Value = 3001
$ws.Range("L22").Value = 7250
$ws.Range("M22").Value = -2706
$ws.Range("N22").Value = -7840
$ws.Range("H27").Value = 4417.3335
$ws.Range("I27").Value = 3001
$ws.Range("J27").Value = 7250
$ws.Range("K27").Value = 3001
$ws.Range("L27").Value = 7250
$ws.Range("M27").Value = -2894
$ws.Range("N27").Value = -7464
$ws.Range("H46").Value = 10265.823
$ws.Range("I46").Value = 4991
$ws.Range("J46").Value = 10595.5
$ws.Range("K46").Value = 4991
$ws.Range("L46").Value = 10595.5
$ws.Range("M46").Value = -4803
$ws.Range("N46").Value = -10971.5
$ws.Range("H61").Value = 1452.4286
$ws.Range("I61").Value = 1294.579
$ws.Range("J61").Value = 2952
$ws.Range("K61").Value = 1294.579
$ws.Range("L61").Value = 2952
$ws.Range("M61").Value = -1092.579
$ws.Range("N61").Value = -3356
$ws.Range("H113").Value = 1452.4286
$ws.Range("I113").Value = 1294.579
$ws.Range("J113").Value = 2952
$ws.Range("K113").Value = 1294.579
$ws.Range("L113").Value = 2952
$ws.Range("M113").Value = 875.421
$ws.Range("N113").Value = -7292
$ws.Range("H126").Value = 6518.5
$ws.Range("I126").Value = 6577.5454
$ws.Range("J126").Value = 6425.7144
$ws.Range("K126").Value = 19732.6362
$ws.Range("L126").Value = 19277.1432
$ws.Range("M126").Value = -17262.6362
$ws.Range("N126").Value = -24217.1432
$ws.Range("H132").Value = 4766.593
$ws.Range("I132").Value = 4519.9565
$ws.Range("K132").Value = 13559.8695
$ws.Range("M132").Value = -11029.8695

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6089.5557
$ws.Range("I126").Value = 6029.5713
$ws.Range("K126").Value = 18088.7139
$ws.Range("M126").Value = -15618.7139
$ws.Range("H132").Value = 2759
$ws.Range("I132").Value = 2852.8572
$ws.Range("J132").Value = 2102
$ws.Range("K132").Value = 8558.571599999999
$ws.Range("L132").Value = 6306
$ws.Range("M132").Value = -6028.571599999999
$ws.Range("N132").Value = -11366

